# Update the LR-pairs data with newly computed TPM-based values.
# The sheet keeps the same 20 columns (A:T) but the data rows now cover
# every Sending-cluster x Target-cluster combination (ECs/FAPs/MuSCs x
# ECs/FAPs/MuSCs, plus Resolving-Mac as a sender), growing from 8 data
# rows (rows 2-9) to 12 data rows (rows 2-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("ECs","Tnfsf10","Tnfrsf11b","ECs",3,1,25.70233233333333,77.10699700000001,0.9446552342719222,0.9446552342719222,1,0.3333333333333333,0.06861733333333334,0.205852,0.01654048691795588,0.01654048691795588,1.763625505160445,15.872629546444,0.01562505754445328,0.01562505754445328),
    @("ECs","Tnfsf10","Tnfrsf11b","FAPs",3,1,25.70233233333333,77.10699700000001,0.9446552342719222,0.9446552342719222,3,1,3.776574666666666,11.329724,0.9103586635352137,0.9103586635352137,97.06677716431422,873.6009944788279,0.8599750765733313,0.8599750765733313),
    @("ECs","Tnfsf10","Tnfrsf11b","MuSCs",3,1,25.70233233333333,77.10699700000001,0.9446552342719222,0.9446552342719222,2,0.6666666666666666,0.303255,0.909765,0.07310084954683041,0.07310084954683042,7.794360791745,70.149247125705,0.06905510015413761,0.06905510015413763),
    @("FAPs","Tnfsf10","Tnfrsf11b","ECs",1,0.3333333333333333,0.234533,0.703599,0.00861995025144722,0.00861995025144722,1,0.3333333333333333,0.06861733333333334,0.205852,0.01654048691795588,0.01654048691795588,0.01609302903866667,0.144837261348,0.0001425781743674933,0.0001425781743674933),
    @("FAPs","Tnfsf10","Tnfrsf11b","FAPs",1,0.3333333333333333,0.234533,0.703599,0.00861995025144722,0.00861995025144722,3,1,3.776574666666666,11.329724,0.9103586635352137,0.9103586635352137,0.8857313862973333,7.971582476675999,0.00784724639064752,0.00784724639064752),
    @("FAPs","Tnfsf10","Tnfrsf11b","MuSCs",1,0.3333333333333333,0.234533,0.703599,0.00861995025144722,0.00861995025144722,2,0.6666666666666666,0.303255,0.909765,0.07310084954683041,0.07310084954683042,0.07112330491499999,0.640109744235,0.0006301256864322061,0.0006301256864322064),
    @("MuSCs","Tnfsf10","Tnfrsf11b","ECs",2,0.6666666666666666,0.2552456666666666,0.765737,0.009381216922838777,0.009381216922838777,1,0.3333333333333333,0.06861733333333334,0.205852,0.01654048691795588,0.01654048691795588,0.01751427699155556,0.157628492924,0.0001551698957867211,0.0001551698957867211),
    @("MuSCs","Tnfsf10","Tnfrsf11b","FAPs",2,0.6666666666666666,0.2552456666666666,0.765737,0.009381216922838777,0.009381216922838777,3,1,3.776574666666666,11.329724,0.9103586635352137,0.9103586635352137,0.9639543185097776,8.675588866587999,0.008540272100209439,0.008540272100209439),
    @("MuSCs","Tnfsf10","Tnfrsf11b","MuSCs",2,0.6666666666666666,0.2552456666666666,0.765737,0.009381216922838777,0.009381216922838777,2,0.6666666666666666,0.303255,0.909765,0.07310084954683041,0.07310084954683042,0.077404524645,0.696640721805,0.0006857749268426168,0.000685774926842617),
    @("Resolving-Mac","Tnfsf10","Tnfrsf11b","ECs",3,1,1.016050666666667,3.048152,0.03734359855379179,0.03734359855379179,1,0.3333333333333333,0.06861733333333334,0.205852,0.01654048691795588,0.01654048691795588,0.06971868727822222,0.627468185504,0.0006176813033483894,0.0006176813033483894),
    @("Resolving-Mac","Tnfsf10","Tnfrsf11b","FAPs",3,1,1.016050666666667,3.048152,0.03734359855379179,0.03734359855379179,3,1,3.776574666666666,11.329724,0.9103586635352137,0.9103586635352137,3.837191207783111,34.534720870048,0.03399606847102543,0.03399606847102543),
    @("Resolving-Mac","Tnfsf10","Tnfrsf11b","MuSCs",3,1,1.016050666666667,3.048152,0.03734359855379179,0.03734359855379179,2,0.6666666666666666,0.303255,0.909765,0.07310084954683041,0.07310084954683042,0.30812244492,2.77310200428,0.002729848779417967,0.002729848779417968)
)

$startRow = 2
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowData = $newRows[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowData[$c]
    }
}
